$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Header row 1: rename C1 from "Дисплей" to the more specific
#    "Дисплей, дюйм" (a brand new shared string is created for this).
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Дисплей, дюйм"

# ---------------------------------------------------------------------------
# 2) sheetView changes: zoom 125 -> 115, new topLeftCell anchor, and the
#    active selection moves from M4 to M10.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 115
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M10").Select()

# ---------------------------------------------------------------------------
# 3) Row 7 becomes a copy of the header row (row 1), reusing its style (s=8)
#    and a taller row height (81pt), with the updated "Дисплей, дюйм" label.
# ---------------------------------------------------------------------------
$ws.Range("B1:M1").Copy()
$ws.Range("B7:M7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B7").Value = "Вес, г"
$ws.Range("C7").Value = "Дисплей, дюйм"
$ws.Range("D7").Value = "Частота обновления экрана, Гц"
$ws.Range("E7").Value = "Разрешение основной камеры, МП"
$ws.Range("F7").Value = "Разрешение фронтальной камеры, МП"
$ws.Range("G7").Value = "Встроенная память, Гб"
$ws.Range("H7").Value = "Оперативная память, Гб"
$ws.Range("I7").Value = "Количество ядер процессора"
$ws.Range("J7").Value = "Частота процессора, ГГц"
$ws.Range("K7").Value = "Емкость аккумулятора мА*ч"
$ws.Range("L7").Value = "Цена, руб."
$ws.Range("M7").Value = "Срок службы, лет"

$ws.Rows.Item(7).RowHeight = 81

# ---------------------------------------------------------------------------
# 4) Rows 8-10 and 12-14: the "Модель" labels in column A switch from the old
#    sample camera names to the real phone names used elsewhere in the sheet,
#    reusing the data-row style (s=4) taken from A2:A4.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A3").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A4").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A8").Value = "Xiaomi POCO M5"
$ws.Range("A9").Value = "Redmi Note 12"
$ws.Range("A10").Value = "POCO X5 Pro"

$ws.Range("A12").Value = "Xiaomi POCO M5"
$ws.Range("A13").Value = "Redmi Note 12"
$ws.Range("A14").Value = "POCO X5 Pro"

# ---------------------------------------------------------------------------
# 5) Rows 8-10: replace the old formula-driven scores with the recomputed
#    literal scores for the three phones (B:M). Downstream formulas in rows
#    12-14 and columns N/O recalc automatically from these.
# ---------------------------------------------------------------------------
$row8 = @(9.1, 9.870000000000001, 7.5, 4.63, 3.13, 5, 7.5, 10, 7.86, 10, 10, 10)
$row9 = @(9.28, 10, 10, 4.63, 6.84, 10, 10, 10, 10, 10, 4.48, 10)
$row10 = @(10, 10, 10, 10, 10, 10, 10, 10, 8.57, 10, 2.21, 2)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $row8[$i]
    $ws.Range($cols[$i] + "9").Value = $row9[$i]
    $ws.Range($cols[$i] + "10").Value = $row10[$i]
}

$wb.Application.Calculate()
